# The sheet originally listed two MS Dhoni (c)† innings vs Chennai Super
# Kings (rows 2 and 3). The "Oct 25 2020" row was a duplicate/bad entry;
# row 2 should instead hold the "Oct 10 2020" data (currently in row 3),
# and the old row 3 should be removed so the sheet shrinks to A1:K2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches the original t="str" cell type) before
# writing numeric-looking values, so Excel doesn't reinterpret them as
# numbers.
$textCells = @("G2", "H2", "I2", "J2", "K2")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = " Oct 10 2020"
$ws.Range("B2").Value = " Dubai (DSC)"
$ws.Range("C2").Value = "RCB won by 37 runs"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "MS Dhoni (c)†"
$ws.Range("G2").Value = "10"
$ws.Range("H2").Value = "6"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "166.66"

# Remove the now-duplicated row 3 entirely; this shifts nothing else and
# shrinks the sheet's used range down to A1:K2.
$ws.Rows.Item(3).Delete()
